$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2561.7856
$ws.Range("I98").Value = 1821.5834
$ws.Range("J98").Value = 7003
$ws.Range("K98").Value = 1821.5834
$ws.Range("L98").Value = 7003
$ws.Range("M98").Value = -323.5834
$ws.Range("N98").Value = -9999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2561.7856
$ws.Range("I122").Value = 1821.5834
$ws.Range("J122").Value = 7003
$ws.Range("K122").Value = 5464.7502
$ws.Range("L122").Value = 21009
$ws.Range("M122").Value = -3014.7502
$ws.Range("N122").Value = -25909

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 2058428.4
$ws.Range("I129").Value = 465
$ws.Range("J129").Value = 3368041.5
$ws.Range("K129").Value = 1395
$ws.Range("L129").Value = 10104124.5
$ws.Range("M129").Value = 3605
$ws.Range("N129").Value = -10114124.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 16668346
$ws.Range("I137").Value = 1732
$ws.Range("J137").Value = 83334800
$ws.Range("K137").Value = 5196
$ws.Range("L137").Value = 250004400
$ws.Range("M137").Value = -2646
$ws.Range("N137").Value = -250009500

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4310.69
$ws.Range("I138").Value = 1250.2941
$ws.Range("J138").Value = 4937.518
$ws.Range("K138").Value = 3750.8823
$ws.Range("L138").Value = 14812.554
$ws.Range("M138").Value = 1389.1177
$ws.Range("N138").Value = -25092.554

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 741.1429000000001
$ws.Range("I110").Value = 958.3333
$ws.Range("J110").Value = 578.25
$ws.Range("K110").Value = 958.3333
$ws.Range("L110").Value = 578.25
$ws.Range("M110").Value = 1086.6667
$ws.Range("N110").Value = -4668.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2380.5908
$ws.Range("I132").Value = 1674.9231
$ws.Range("J132").Value = 3399.889
$ws.Range("K132").Value = 5024.7693
$ws.Range("L132").Value = 10199.667
$ws.Range("M132").Value = -2494.7693
$ws.Range("N132").Value = -15259.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3890.2222
$ws.Range("I134").Value = 3578
$ws.Range("J134").Value = 4140
$ws.Range("K134").Value = 10734
$ws.Range("L134").Value = 12420
$ws.Range("M134").Value = -8199
$ws.Range("N134").Value = -17490

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 20.625
$ws.Range("I12").Value = 51
$ws.Range("J12").Value = 10.5
$ws.Range("K12").Value = 153
$ws.Range("L12").Value = 31.5
$ws.Range("M12").Value = 20
$ws.Range("N12").Value = -377.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 852.0202
$ws.Range("I68").Value = 621.0333000000001
$ws.Range("J68").Value = 952.4493
$ws.Range("K68").Value = 1863.0999
$ws.Range("L68").Value = 2857.3479
$ws.Range("M68").Value = -1052.0999
$ws.Range("N68").Value = -4479.3479

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 852.0202
$ws.Range("I71").Value = 621.0333000000001
$ws.Range("J71").Value = 952.4493
$ws.Range("K71").Value = 5589.2997
$ws.Range("L71").Value = 8572.0437
$ws.Range("M71").Value = -1533.2997
$ws.Range("N71").Value = -16684.0437

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 21584.725
$ws.Range("I131").Value = 144785.72
$ws.Range("J131").Value = 1984.5682
$ws.Range("K131").Value = 434357.16
$ws.Range("L131").Value = 5953.7046
$ws.Range("M131").Value = -429317.16
$ws.Range("N131").Value = -16033.7046

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4875
$ws.Range("I80").Value = 4875
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 4875
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -3877
$ws.Range("N80").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4875
$ws.Range("I83").Value = 4875
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 24375
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -19383
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6495006
$ws.Range("I122").Value = 12987012
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 38961036
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -38958586
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 39993.332
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 39993.332
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 39993.332
$ws.Range("N124").Value = -49813.332

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H128").Value = 35986.785
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 35986.785
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 35986.785
$ws.Range("N128").Value = -45946.785

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 34120
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 34120
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 34120
$ws.Range("N133").Value = -44240

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 19818.166
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 19818.166
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 19818.166
$ws.Range("N138").Value = -30098.166

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 70750
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 70750
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 70750
$ws.Range("N140").Value = -81110

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 42999.5
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 42999.5
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 42999.5
$ws.Range("N141").Value = -53359.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1281876.8
$ws.Range("I40").Value = 1281876.8
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1281876.8
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1281740.8
$ws.Range("N40").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6714.3477
$ws.Range("I122").Value = 6864.773
$ws.Range("J122").Value = 3405
$ws.Range("K122").Value = 20594.319
$ws.Range("L122").Value = 10215
$ws.Range("M122").Value = -18144.319
$ws.Range("N122").Value = -15115

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3975.9333
$ws.Range("I136").Value = 1323.9
$ws.Range("J136").Value = 9280
$ws.Range("K136").Value = 3971.7
$ws.Range("L136").Value = 27840
$ws.Range("M136").Value = -1421.7
$ws.Range("N136").Value = -32940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 38700
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 38700
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 38700
$ws.Range("N139").Value = -48980

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I126").Value = 850
$ws.Range("J126").Value = 900
$ws.Range("K126").Value = 2550
$ws.Range("L126").Value = 2700
$ws.Range("M126").Value = -80
$ws.Range("N126").Value = -7640

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H140").Value = 56200
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 56200
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 56200
$ws.Range("N140").Value = -66560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 44916.668
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 44916.668
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 44916.668
$ws.Range("N141").Value = -55276.668
